$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently ends at row 294. We need to append 9 new documentation
# rows (295-303). Copy the last existing row 9 times and insert the copies
# right after it so every new row inherits the existing date-style (s="1")
# formatting instead of Excel inventing a brand-new custom number format.
for ($i = 0; $i -lt 9; $i++) {
    $ws.Rows(294).Copy()
    $ws.Rows(295).Insert(-4121)
}

# Row 295: 2018-10-14, photo 1
$ws.Range("A295").Value = 43387
$ws.Range("B295").Value = "2018-10-14 1.JPG"
$ws.Range("C295").Value = "Drohne stürzt bei versuch ab, Kabel schmoren durch das Zimmer stikt. Gott sei dank steht das Haus noch"

# Row 296: 2018-10-14, photo 2
$ws.Range("A296").Value = 43387
$ws.Range("B296").Value = "2018-10-14 2.JPG"
$ws.Range("C296").Value = "Höhensensor ist gekommen zum Ausprobieren"

# Row 297: 2018-10-14, photo 3
$ws.Range("A297").Value = 43387
$ws.Range("B297").Value = "2018-10-14 3.JPG"
$ws.Range("C297").Value = "Bei einem weiteren Flugversuch macht ein Motor seltsame geräusche. Bis man den winzigen stein gefunden hat vergehen auch wieder ewigkeiten"

# Row 298: 2018-10-21, photo 1
$ws.Range("A298").Value = 43394
$ws.Range("B298").Value = "2018-10-21 1.JPG"
$ws.Range("C298").Value = "Upgrade auf GPS-Sensor (unten) und umrüstung auf den RasPi3"

# Row 299: 2018-10-21, photo 2
$ws.Range("A299").Value = 43394
$ws.Range("B299").Value = "2018-10-21 2.JPG"
$ws.Range("C299").Value = "Die Löstspitze ist inzwischen auch durchgerostet und fast zerstört XD"

# Row 300: 2018-10-21, photo 3
$ws.Range("A300").Value = 43394
$ws.Range("B300").Value = "2018-10-21 3.JPG"
$ws.Range("C300").Value = "Die Löstspitze ist inzwischen auch durchgerostet und fast zerstört XD"

# Row 301: 2018-10-30
$ws.Range("A301").Value = 43403
$ws.Range("B301").Value = "2018-10-30.JPG"
$ws.Range("C301").Value = "Wie viel Pech muss man haben, dass der SD-Kartenleser am Raspi kaputt geht??? Nichts, was sich nicht lösen ließe"

# Row 302: 2018-11-17, no image
$ws.Range("A302").Value = 43421
$ws.Range("B302").ClearContents()
$ws.Range("C302").Value = "Vorletzter Entwicklungstag"

# Row 303: 2018-11-18, no image
$ws.Range("A303").Value = 43422
$ws.Range("B303").ClearContents()
$ws.Range("C303").Value = "Letzter Entwicklungstag"

# Match the author's final view/selection state from the diff.
$ws.Application.GoTo($ws.Range("C301"), $true) | Out-Null
$ws.Range("C301").Select() | Out-Null
